$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 351; everything from old row 351 onward shifts
# down by 3 rows (old row 351 -> new row 354, ... old row 400 -> new row 403).
$ws.Rows("351:353").Insert()

# The three new rows carry the same constant columns as the rest of the
# "Comercializadora del Agro de Limari - Limon" dataset, a new reporting
# date (2021-10-07 = serial 44476), and the weekly 1a/2a/3a amarillo prices.
$rows = @(
    @{ Row = 351; L = "1a amarillo"; M = 900; N = 3800; O = 4000; P = 3900; S = 244 },
    @{ Row = 352; L = "2a amarillo"; M = 750; N = 2800; O = 3000; P = 2900; S = 181 },
    @{ Row = 353; L = "3a amarillo"; M = 540; N = 1800; O = 2000; P = 1900; S = 119 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44476
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102003
    $ws.Cells.Item($row, 10).Value = "Limón"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/malla 16 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 16
}
